$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 259, pushing existing data (old rows 259-262) down to 261-264
$ws.Rows.Item(259).Resize(2).Insert()

# Copy date formatting for column D from the row right below (now row 261) into the new rows
$ws.Range("D259:D260").NumberFormat = $ws.Range("D261").NumberFormat

# Row 259 - Black Amber
$ws.Cells.Item(259, 1).Value = 4
$ws.Cells.Item(259, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(259, 3).Value = "Los Lagos"
$ws.Cells.Item(259, 4).Value = 44939
$ws.Cells.Item(259, 5).Value = 10
$ws.Cells.Item(259, 6).Value = "Fruta"
$ws.Cells.Item(259, 7).Value = 100103
$ws.Cells.Item(259, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(259, 9).Value = 100103002
$ws.Cells.Item(259, 10).Value = "Ciruela"
$ws.Cells.Item(259, 11).Value = "Black Amber"
$ws.Cells.Item(259, 12).Value = "Primera"
$ws.Cells.Item(259, 13).Value = 600
$ws.Cells.Item(259, 14).Value = 17000
$ws.Cells.Item(259, 15).Value = 18000
$ws.Cells.Item(259, 16).Value = 17500
$ws.Cells.Item(259, 17).Value = "$/caja 14 kilos granel"
$ws.Cells.Item(259, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(259, 19).Value = 1250
$ws.Cells.Item(259, 20).Value = 14

# Row 260 - Lemon
$ws.Cells.Item(260, 1).Value = 4
$ws.Cells.Item(260, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(260, 3).Value = "Los Lagos"
$ws.Cells.Item(260, 4).Value = 44939
$ws.Cells.Item(260, 5).Value = 10
$ws.Cells.Item(260, 6).Value = "Fruta"
$ws.Cells.Item(260, 7).Value = 100103
$ws.Cells.Item(260, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(260, 9).Value = 100103002
$ws.Cells.Item(260, 10).Value = "Ciruela"
$ws.Cells.Item(260, 11).Value = "Lemon"
$ws.Cells.Item(260, 12).Value = "Primera"
$ws.Cells.Item(260, 13).Value = 600
$ws.Cells.Item(260, 14).Value = 17000
$ws.Cells.Item(260, 15).Value = 18000
$ws.Cells.Item(260, 16).Value = 17500
$ws.Cells.Item(260, 17).Value = "$/caja 14 kilos granel"
$ws.Cells.Item(260, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(260, 19).Value = 1250
$ws.Cells.Item(260, 20).Value = 14
